$d = $word.ActiveDocument

$d.Content.Find.Execute("989×6=", $true, $false, $false, $false, $false, $true, 1, $false, "141×2=", 2) | Out-Null
$d.Content.Find.Execute("488×6=", $true, $false, $false, $false, $false, $true, 1, $false, "756×9=", 2) | Out-Null
$d.Content.Find.Execute("127×5=", $true, $false, $false, $false, $false, $true, 1, $false, "269×4=", 2) | Out-Null
$d.Content.Find.Execute("918×3=", $true, $false, $false, $false, $false, $true, 1, $false, "287×5=", 2) | Out-Null
$d.Content.Find.Execute("520×9=", $true, $false, $false, $false, $false, $true, 1, $false, "493×9=", 2) | Out-Null
$d.Content.Find.Execute("487×6=", $true, $false, $false, $false, $false, $true, 1, $false, "915×6=", 2) | Out-Null
$d.Content.Find.Execute("684×9=", $true, $false, $false, $false, $false, $true, 1, $false, "302×4=", 2) | Out-Null
$d.Content.Find.Execute("541×2=", $true, $false, $false, $false, $false, $true, 1, $false, "454×9=", 2) | Out-Null
$d.Content.Find.Execute("354×4=", $true, $false, $false, $false, $false, $true, 1, $false, "602×9=", 2) | Out-Null
$d.Content.Find.Execute("540×4=", $true, $false, $false, $false, $false, $true, 1, $false, "933×9=", 2) | Out-Null
$d.Content.Find.Execute("606×7=", $true, $false, $false, $false, $false, $true, 1, $false, "401×9=", 2) | Out-Null
$d.Content.Find.Execute("509×2=", $true, $false, $false, $false, $false, $true, 1, $false, "634×8=", 2) | Out-Null
$d.Content.Find.Execute("836×2=", $true, $false, $false, $false, $false, $true, 1, $false, "941×4=", 2) | Out-Null
$d.Content.Find.Execute("153×4=", $true, $false, $false, $false, $false, $true, 1, $false, "526×9=", 2) | Out-Null
$d.Content.Find.Execute("957×7=", $true, $false, $false, $false, $false, $true, 1, $false, "877×5=", 2) | Out-Null
$d.Content.Find.Execute("274×6=", $true, $false, $false, $false, $false, $true, 1, $false, "916×8=", 2) | Out-Null
$d.Content.Find.Execute("156×3=", $true, $false, $false, $false, $false, $true, 1, $false, "979×6=", 2) | Out-Null
$d.Content.Find.Execute("575×6=", $true, $false, $false, $false, $false, $true, 1, $false, "654×3=", 2) | Out-Null
$d.Content.Find.Execute("295×3=", $true, $false, $false, $false, $false, $true, 1, $false, "223×3=", 2) | Out-Null
$d.Content.Find.Execute("162×3=", $true, $false, $false, $false, $false, $true, 1, $false, "987×6=", 2) | Out-Null
$d.Content.Find.Execute("985×6=", $true, $false, $false, $false, $false, $true, 1, $false, "454×9=", 2) | Out-Null
$d.Content.Find.Execute("867×2=", $true, $false, $false, $false, $false, $true, 1, $false, "830×8=", 2) | Out-Null
$d.Content.Find.Execute("532×3=", $true, $false, $false, $false, $false, $true, 1, $false, "649×6=", 2) | Out-Null
$d.Content.Find.Execute("759×3=", $true, $false, $false, $false, $false, $true, 1, $false, "177×9=", 2) | Out-Null
$d.Content.Find.Execute("142×6=", $true, $false, $false, $false, $false, $true, 1, $false, "958×8=", 2) | Out-Null
